$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 145
$ws.Cells.Item(12, 9).Value = 90
$ws.Cells.Item(12, 11).Value = 90
$ws.Cells.Item(12, 13).Value = 80
$ws.Cells.Item(32, 8).Value = 854.8889
$ws.Cells.Item(32, 9).Value = 849.875
$ws.Cells.Item(32, 11).Value = 849.875
$ws.Cells.Item(32, 13).Value = -523.875
$ws.Cells.Item(34, 8).Value = 10316.333
$ws.Cells.Item(34, 9).Value = 7999.5
$ws.Cells.Item(34, 10).Value = 14950
$ws.Cells.Item(34, 11).Value = 7999.5
$ws.Cells.Item(34, 12).Value = 14950
$ws.Cells.Item(34, 13).Value = -7796.5
$ws.Cells.Item(34, 14).Value = -15356
$ws.Cells.Item(36, 8).Value = 10316.333
$ws.Cells.Item(36, 9).Value = 7999.5
$ws.Cells.Item(36, 10).Value = 14950
$ws.Cells.Item(36, 11).Value = 7999.5
$ws.Cells.Item(36, 12).Value = 14950
$ws.Cells.Item(36, 13).Value = -7284.5
$ws.Cells.Item(36, 14).Value = -16380
$ws.Cells.Item(113, 8).Value = 1995
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(138, 8).Value = 2639.4722
$ws.Cells.Item(138, 10).Value = 5426.7144
$ws.Cells.Item(138, 12).Value = 16280.1432
$ws.Cells.Item(138, 14).Value = -26560.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 10124.75
$ws.Cells.Item(10, 10).Value = 10333
$ws.Cells.Item(10, 12).Value = 10333
$ws.Cells.Item(10, 14).Value = -10673
$ws.Cells.Item(32, 8).Value = 3164.1719
$ws.Cells.Item(32, 9).Value = 2360.5344
$ws.Cells.Item(32, 11).Value = 2360.5344
$ws.Cells.Item(32, 13).Value = -2073.5344
$ws.Cells.Item(33, 8).Value = 1513
$ws.Cells.Item(33, 9).Value = 1513
$ws.Cells.Item(33, 11).Value = 1513
$ws.Cells.Item(33, 13).Value = -1184
$ws.Cells.Item(36, 8).Value = 5310.4
$ws.Cells.Item(36, 9).Value = 5310.4
$ws.Cells.Item(36, 11).Value = 5310.4
$ws.Cells.Item(36, 13).Value = -4964.4
$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 13).ClearContents()
$ws.Cells.Item(123, 8).Value = 80000
$ws.Cells.Item(123, 10).Value = 80000
$ws.Cells.Item(123, 12).Value = 80000
$ws.Cells.Item(123, 14).Value = -89800

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 88.666664
$ws.Cells.Item(22, 9).Value = 88.666664
$ws.Cells.Item(22, 11).Value = 88.666664
$ws.Cells.Item(22, 13).Value = 84.333336

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 160.11539
$ws.Cells.Item(7, 9).Value = 123.7619
$ws.Cells.Item(7, 11).Value = 123.7619
$ws.Cells.Item(7, 13).Value = -10.7619
$ws.Cells.Item(29, 8).Value = 8118.5
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 8118.5
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 8118.5
$ws.Cells.Item(29, 13).ClearContents()
$ws.Cells.Item(29, 14).Value = -8704.5
$ws.Cells.Item(35, 8).Value = 479.41666
$ws.Cells.Item(60, 8).Value = 13409
$ws.Cells.Item(60, 10).Value = 32494
$ws.Cells.Item(60, 12).Value = 32494
$ws.Cells.Item(60, 14).Value = -33516
$ws.Cells.Item(68, 8).Value = 37999.668
$ws.Cells.Item(68, 10).Value = 37999.668
$ws.Cells.Item(68, 12).Value = 37999.668
$ws.Cells.Item(68, 14).Value = -39497.668
$ws.Cells.Item(71, 8).Value = 37999.668
$ws.Cells.Item(71, 10).Value = 37999.668
$ws.Cells.Item(71, 12).Value = 113999.004
$ws.Cells.Item(71, 14).Value = -121487.004
$ws.Cells.Item(74, 8).Value = 47567
$ws.Cells.Item(74, 10).Value = 47567
$ws.Cells.Item(74, 12).Value = 47567
$ws.Cells.Item(74, 14).Value = -49315
$ws.Cells.Item(77, 8).Value = 47567
$ws.Cells.Item(77, 10).Value = 47567
$ws.Cells.Item(77, 12).Value = 142701
$ws.Cells.Item(77, 14).Value = -151437

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 112157.336
$ws.Cells.Item(2, 9).Value = 200076
$ws.Cells.Item(2, 10).Value = 2259
$ws.Cells.Item(2, 11).Value = 1200456
$ws.Cells.Item(2, 12).Value = 13554
$ws.Cells.Item(2, 13).Value = -1200343
$ws.Cells.Item(2, 14).Value = -13780
$ws.Cells.Item(11, 8).Value = 342.66666
$ws.Cells.Item(11, 10).Value = 999
$ws.Cells.Item(11, 12).Value = 2997
$ws.Cells.Item(11, 14).Value = -3277
$ws.Cells.Item(26, 8).Value = 750
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 750
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 2250
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(26, 14).Value = -2826
$ws.Cells.Item(34, 8).Value = 1759.75
$ws.Cells.Item(34, 10).Value = 2500
$ws.Cells.Item(34, 12).Value = 7500
$ws.Cells.Item(34, 14).Value = -7668
$ws.Cells.Item(39, 8).Value = 1121.5454
$ws.Cells.Item(39, 10).Value = 1739.6
$ws.Cells.Item(39, 12).Value = 5218.799999999999
$ws.Cells.Item(39, 14).Value = -5806.799999999999
$ws.Cells.Item(55, 8).Value = 93481.63
$ws.Cells.Item(55, 10).Value = 3916.5
$ws.Cells.Item(55, 12).Value = 11749.5
$ws.Cells.Item(55, 14).Value = -12103.5
$ws.Cells.Item(120, 8).Value = 8098.6665
$ws.Cells.Item(120, 10).Value = 9285.714
$ws.Cells.Item(120, 12).Value = 27857.142
$ws.Cells.Item(120, 14).Value = -37533.142

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 1644.5
$ws.Cells.Item(21, 9).Value = 1644.5
$ws.Cells.Item(21, 11).Value = 1644.5
$ws.Cells.Item(21, 13).Value = -1471.5
$ws.Cells.Item(30, 8).Value = 1644.5
$ws.Cells.Item(30, 9).Value = 1644.5
$ws.Cells.Item(30, 11).Value = 1644.5
$ws.Cells.Item(30, 13).Value = -1539.5
$ws.Cells.Item(113, 8).Value = 99999
$ws.Cells.Item(113, 9).Value = 99999
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 99999
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(113, 13).Value = -97829
$ws.Cells.Item(125, 8).Value = 81999
$ws.Cells.Item(125, 10).Value = 81999
$ws.Cells.Item(125, 12).Value = 81999
$ws.Cells.Item(125, 14).Value = -86919

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 5875
$ws.Cells.Item(4, 9).Value = 9000
$ws.Cells.Item(4, 11).Value = 9000
$ws.Cells.Item(4, 13).Value = -8887
$ws.Cells.Item(22, 8).Value = 6018.3335
$ws.Cells.Item(22, 10).Value = 6018.3335
$ws.Cells.Item(22, 12).Value = 6018.3335
$ws.Cells.Item(22, 14).Value = -6608.3335
$ws.Cells.Item(23, 8).Value = 10502000
$ws.Cells.Item(23, 9).Value = 10502000
$ws.Cells.Item(23, 11).Value = 10502000
$ws.Cells.Item(23, 13).Value = -10501770
$ws.Cells.Item(27, 8).Value = 6018.3335
$ws.Cells.Item(27, 10).Value = 6018.3335
$ws.Cells.Item(27, 12).Value = 6018.3335
$ws.Cells.Item(27, 14).Value = -6232.3335
$ws.Cells.Item(28, 8).Value = 5875
$ws.Cells.Item(28, 9).Value = 9000
$ws.Cells.Item(28, 11).Value = 9000
$ws.Cells.Item(28, 13).Value = -8768
$ws.Cells.Item(37, 8).Value = 5875
$ws.Cells.Item(37, 9).Value = 9000
$ws.Cells.Item(37, 11).Value = 9000
$ws.Cells.Item(37, 13).Value = -8893
$ws.Cells.Item(41, 8).Value = 27000
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 27000
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 27000
$ws.Cells.Item(41, 13).ClearContents()
$ws.Cells.Item(41, 14).Value = -27876
$ws.Cells.Item(43, 8).Value = 392469.06
$ws.Cells.Item(43, 9).Value = 7999
$ws.Cells.Item(43, 10).Value = 424508.25
$ws.Cells.Item(43, 11).Value = 7999
$ws.Cells.Item(43, 12).Value = 424508.25
$ws.Cells.Item(43, 13).Value = -7806
$ws.Cells.Item(43, 14).Value = -424894.25
$ws.Cells.Item(47, 8).Value = 24000
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 13).ClearContents()
$ws.Cells.Item(52, 8).Value = 24000
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 20005
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 20005
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 20005
$ws.Cells.Item(11, 13).ClearContents()
$ws.Cells.Item(11, 14).Value = -20289
$ws.Cells.Item(18, 8).Value = 13502.5
$ws.Cells.Item(18, 10).Value = 13502.5
$ws.Cells.Item(18, 12).Value = 13502.5
$ws.Cells.Item(18, 14).Value = -13848.5
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 13).ClearContents()
$ws.Cells.Item(81, 8).Value = 2960.8
$ws.Cells.Item(81, 9).Value = 3212.5557
$ws.Cells.Item(81, 10).Value = 695
$ws.Cells.Item(81, 11).Value = 6425.1114
$ws.Cells.Item(81, 12).Value = 1390
$ws.Cells.Item(81, 13).Value = -5364.1114
$ws.Cells.Item(81, 14).Value = -3512
$ws.Cells.Item(84, 8).Value = 2960.8
$ws.Cells.Item(84, 9).Value = 3212.5557
$ws.Cells.Item(84, 10).Value = 695
$ws.Cells.Item(84, 11).Value = 32125.557
$ws.Cells.Item(84, 12).Value = 6950
$ws.Cells.Item(84, 13).Value = -26821.557
$ws.Cells.Item(84, 14).Value = -17558
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).ClearContents()
